$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells originally hold plain text (inline strings) such as
# "67.848.39", "1.00", "0.000274", etc. Excel auto-converts plausible
# numeric-looking strings to actual numbers when assigned via .Value,
# which would corrupt values like "1.00" -> 1 or "0.524" -> 0.524000...
# Force the Text number format first so the assignment is kept verbatim
# as a string, matching the original inline-string cell type.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.848.39"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.849.24"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.50%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.19"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.80%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.08"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +1.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.848.85"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -1.48%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.12%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.164"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -0.46%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.32"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.78%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.455"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000247"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.75"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.499.82"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.869.66"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.914.41"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -1.06%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.08"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.39%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.34"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.90"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "462.93"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.728"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +1.55%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -3.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.15"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -1.36%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.08"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.96%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.97"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.94"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +0.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.003.11"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -1.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.73"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.31"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -2.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.97"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -2.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.829.38"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.53%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -2.29%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -1.70%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.139"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.18%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.88"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +0.30%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.27"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +7.21%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.53%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "426.78"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.17%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.97"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "47.13"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.48"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.79%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "26.73"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.37%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "FLOKI"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000274"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +3.58%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "143.54"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "40.19"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.99%  "
